# This workbook contains a weekly price report table for
# "Fruta, Feria Lagunitas de Puerto Montt - Kiwi".
# The commit adds one more week of data: two new report rows
# (quality "Especial" and "Primera") are inserted right before the
# existing row 43, pushing the rest of the table down by two rows
# (old row 43 -> new row 45, ..., old row 146 -> new row 148).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (row 43),
# which shifts every following row down by two and extends the
# used range from A1:T146 to A1:T148.
$ws.Rows("43:44").Insert()

# --- Fill in the first new row (row 43) ---
$ws.Cells.Item(43, 1).Value = 4
$ws.Cells.Item(43, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(43, 3).Value = "Los Lagos"
$ws.Cells.Item(43, 4).Value = 44498
$ws.Cells.Item(43, 5).Value = 10
$ws.Cells.Item(43, 6).Value = "Fruta"
$ws.Cells.Item(43, 7).Value = 100101
$ws.Cells.Item(43, 8).Value = "Berries"
$ws.Cells.Item(43, 9).Value = 100101007
$ws.Cells.Item(43, 10).Value = "Kiwi"
$ws.Cells.Item(43, 11).Value = "Hayward"
$ws.Cells.Item(43, 12).Value = "Especial"
$ws.Cells.Item(43, 13).Value = 200
$ws.Cells.Item(43, 14).Value = 23000
$ws.Cells.Item(43, 15).Value = 23000
$ws.Cells.Item(43, 16).Value = 23000
$ws.Cells.Item(43, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(43, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(43, 19).Value = 1533
$ws.Cells.Item(43, 20).Value = 15

# --- Fill in the second new row (row 44) ---
$ws.Cells.Item(44, 1).Value = 4
$ws.Cells.Item(44, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(44, 3).Value = "Los Lagos"
$ws.Cells.Item(44, 4).Value = 44498
$ws.Cells.Item(44, 5).Value = 10
$ws.Cells.Item(44, 6).Value = "Fruta"
$ws.Cells.Item(44, 7).Value = 100101
$ws.Cells.Item(44, 8).Value = "Berries"
$ws.Cells.Item(44, 9).Value = 100101007
$ws.Cells.Item(44, 10).Value = "Kiwi"
$ws.Cells.Item(44, 11).Value = "Hayward"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 600
$ws.Cells.Item(44, 14).Value = 16000
$ws.Cells.Item(44, 15).Value = 17000
$ws.Cells.Item(44, 16).Value = 16500
$ws.Cells.Item(44, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(44, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(44, 19).Value = 1100
$ws.Cells.Item(44, 20).Value = 15

# Make sure the number format applied to the new date cells (D43, D44)
# matches the date format used throughout column D (style index 2).
$dateFmt = $ws.Cells.Item(45, 4).NumberFormat
$ws.Cells.Item(43, 4).NumberFormat = $dateFmt
$ws.Cells.Item(44, 4).NumberFormat = $dateFmt
